# Update "want to go" counts (column F) on several sheets to reflect the
# refreshed scrape output (gh-pages rebuild at commit 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 9
$ws1.Range("F3").Value = 2719
$ws1.Range("F5").Value = 19527
$ws1.Range("F10").Value = 429
$ws1.Range("F22").Value = 96

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 191
$ws2.Range("F4").Value = 9
$ws2.Range("F12").Value = 10

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 5995
$ws3.Range("F4").Value = 586

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 5995
$ws4.Range("F4").Value = 586
$ws4.Range("F5").Value = 191
$ws4.Range("F6").Value = 9
$ws4.Range("F8").Value = 2719
$ws4.Range("F10").Value = 19528
$ws4.Range("F11").Value = 9
$ws4.Range("F20").Value = 429
$ws4.Range("F30").Value = 10
$ws4.Range("F47").Value = 96
